$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 7145382
$ws.Range("J17").Value = 7145382
$ws.Range("L17").Value = 21436146
$ws.Range("N17").Value = -21436482

$ws.Range("H19").Value = 2393.8845
$ws.Range("I19").Value = 1377.0667
$ws.Range("K19").Value = 1377.0667
$ws.Range("M19").Value = -1202.0667

$ws.Range("H39").Value = 799.1667
$ws.Range("I39").Value = 600
$ws.Range("J39").Value = 898.75
$ws.Range("K39").Value = 1800
$ws.Range("L39").Value = 2696.25
$ws.Range("M39").Value = -1504
$ws.Range("N39").Value = -3288.25

$ws.Range("H58").Value = 1827.8572
$ws.Range("I58").Value = 298.66666
$ws.Range("J58").Value = 2974.75
$ws.Range("K58").Value = 895.9999799999999
$ws.Range("L58").Value = 8924.25
$ws.Range("M58").Value = -745.9999799999999
$ws.Range("N58").Value = -9224.25

$ws.Range("H64").Value = 11858.333
$ws.Range("I64").Value = 4160
$ws.Range("J64").Value = 13398
$ws.Range("K64").Value = 4160
$ws.Range("L64").Value = 13398
$ws.Range("M64").Value = -3912
$ws.Range("N64").Value = -13894

$ws.Range("H67").Value = 11858.333
$ws.Range("I67").Value = 4160
$ws.Range("J67").Value = 13398
$ws.Range("K67").Value = 4160
$ws.Range("L67").Value = 13398
$ws.Range("M67").Value = -3302
$ws.Range("N67").Value = -15114

$ws.Range("H80").Value = 10415.9
$ws.Range("I80").Value = 10415.3
$ws.Range("J80").Value = 10416.5
$ws.Range("K80").Value = 31245.9
$ws.Range("L80").Value = 31249.5
$ws.Range("M80").Value = -30247.9
$ws.Range("N80").Value = -33245.5

$ws.Range("H83").Value = 10415.9
$ws.Range("I83").Value = 10415.3
$ws.Range("J83").Value = 10416.5
$ws.Range("K83").Value = 93737.7
$ws.Range("L83").Value = 93748.5
$ws.Range("M83").Value = -88745.7
$ws.Range("N83").Value = -103732.5

$ws.Range("H106").Value = 3405.75
$ws.Range("I106").Value = 3443.0908
$ws.Range("K106").Value = 3443.0908
$ws.Range("M106").Value = -2812.0908

$ws.Range("H138").Value = 5099.5557
$ws.Range("I138").Value = 13623.75
$ws.Range("J138").Value = 3617.087
$ws.Range("K138").Value = 40871.25
$ws.Range("L138").Value = 10851.261
$ws.Range("M138").Value = -35731.25
$ws.Range("N138").Value = -21131.261

$ws.Range("H141").Value = 1971.3684
$ws.Range("I141").Value = 1971.3684
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 5914.1052
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -734.1052
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10053.796
$ws.Range("I32").Value = 10056.357
$ws.Range("K32").Value = 10056.357
$ws.Range("M32").Value = -9769.357

$ws.Range("H74").Value = 2912.2964
$ws.Range("I74").Value = 1922.4
$ws.Range("K74").Value = 1922.4
$ws.Range("M74").Value = -1048.4

$ws.Range("H77").Value = 2912.2964
$ws.Range("I77").Value = 1922.4
$ws.Range("K77").Value = 9612
$ws.Range("M77").Value = -5244

$ws.Range("H92").Value = 65400
$ws.Range("J92").Value = 65400
$ws.Range("L92").Value = 65400
$ws.Range("N92").Value = -70392

$ws.Range("H102").Value = 3329
$ws.Range("I102").Value = 2308.8823
$ws.Range("K102").Value = 2308.8823
$ws.Range("M102").Value = -686.8823000000002

$ws.Range("H109").Value = 78788
$ws.Range("J109").Value = 78788
$ws.Range("L109").Value = 78788
$ws.Range("N109").Value = -81562

$ws.Range("H132").Value = 1067.7561
$ws.Range("I132").Value = 963.5278
$ws.Range("K132").Value = 2890.5834
$ws.Range("M132").Value = -360.5834

$ws.Range("H135").Value = 77554
$ws.Range("J135").Value = 77554
$ws.Range("L135").Value = 77554
$ws.Range("N135").Value = -87694

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7943.6665
$ws.Range("I20").Value = 6998.25
$ws.Range("K20").Value = 6998.25
$ws.Range("M20").Value = -6751.25

$ws.Range("H37").Value = 11931.25
$ws.Range("I37").Value = 13512.5
$ws.Range("K37").Value = 13512.5
$ws.Range("M37").Value = -13375.5

$ws.Range("H99").Value = 5105.222
$ws.Range("I99").Value = 4550
$ws.Range("K99").Value = 4550
$ws.Range("M99").Value = -3052

$ws.Range("H134").Value = 2345.3333
$ws.Range("J134").Value = 2833
$ws.Range("L134").Value = 8499
$ws.Range("N134").Value = -13569

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 58.75
$ws.Range("I7").Value = 58.75
$ws.Range("K7").Value = 58.75
$ws.Range("M7").Value = 54.25

$ws.Range("H64").Value = 51326.168
$ws.Range("J64").Value = 51326.168
$ws.Range("L64").Value = 51326.168
$ws.Range("N64").Value = -51822.168

$ws.Range("H67").Value = 51326.168
$ws.Range("J67").Value = 51326.168
$ws.Range("L67").Value = 51326.168
$ws.Range("N67").Value = -53042.168

$ws.Range("H134").Value = 3709.818
$ws.Range("I134").Value = 2447.625
$ws.Range("J134").Value = 7075.6665
$ws.Range("K134").Value = 7342.875
$ws.Range("L134").Value = 21226.9995
$ws.Range("M134").Value = -4807.875
$ws.Range("N134").Value = -26296.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 86250
$ws.Range("I56").Value = 86250
$ws.Range("K56").Value = 86250
$ws.Range("M56").Value = -85720

$ws.Range("H137").Value = 1938.55
$ws.Range("I137").Value = 1595.0834
$ws.Range("J137").Value = 2453.75
$ws.Range("K137").Value = 4785.2502
$ws.Range("L137").Value = 7361.25
$ws.Range("M137").Value = 314.7497999999996
$ws.Range("N137").Value = -17561.25

$ws.Range("H140").Value = 2073.3333
$ws.Range("I140").Value = 2077.7058
$ws.Range("K140").Value = 6233.117400000001
$ws.Range("M140").Value = -1053.117400000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 39
$ws.Range("I2").Value = 37.791668
$ws.Range("J2").Value = 46.25
$ws.Range("K2").Value = 37.791668
$ws.Range("L2").Value = 46.25
$ws.Range("M2").Value = 75.208332
$ws.Range("N2").Value = -272.25

$ws.Range("H17").Value = 9003.333000000001
$ws.Range("I17").Value = 223
$ws.Range("J17").Value = 52905
$ws.Range("K17").Value = 223
$ws.Range("L17").Value = 52905
$ws.Range("M17").Value = -55
$ws.Range("N17").Value = -53241

$ws.Range("H132").Value = 3711.5625
$ws.Range("I132").Value = 2299.125
$ws.Range("J132").Value = 5124
$ws.Range("K132").Value = 6897.375
$ws.Range("L132").Value = 15372
$ws.Range("M132").Value = -4367.375
$ws.Range("N132").Value = -20432

$ws.Range("H134").Value = 54337.848
$ws.Range("J134").Value = 54337.848
$ws.Range("L134").Value = 163013.544
$ws.Range("N134").Value = -168083.544

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 311.23077
$ws.Range("I16").Value = 332.36365
$ws.Range("J16").Value = 195
$ws.Range("K16").Value = 332.36365
$ws.Range("L16").Value = 195
$ws.Range("M16").Value = -162.36365
$ws.Range("N16").Value = -535

$ws.Range("H136").Value = 54058660
$ws.Range("I136").Value = 43483430
$ws.Range("J136").Value = 71432250
$ws.Range("K136").Value = 130450290
$ws.Range("L136").Value = 214296750
$ws.Range("M136").Value = -130447740
$ws.Range("N136").Value = -214301850

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2327.4614
$ws.Range("I81").Value = 2568.9092
$ws.Range("K81").Value = 5137.8184
$ws.Range("M81").Value = -4076.8184

$ws.Range("H84").Value = 2327.4614
$ws.Range("I84").Value = 2568.9092
$ws.Range("K84").Value = 25689.092
$ws.Range("M84").Value = -20385.092

$ws.Range("H136").Value = 5824.3887
$ws.Range("I136").Value = 4530.4595
$ws.Range("J136").Value = 8640.588
$ws.Range("K136").Value = 13591.3785
$ws.Range("L136").Value = 25921.764
$ws.Range("M136").Value = -11041.3785
$ws.Range("N136").Value = -31021.764
